$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.355.63'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '1.889.16'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.12'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4811'
$ws.Range("E7").Value = '  -2.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2897'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06606'
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("D10").Value = '1.906.76'
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.93'
$ws.Range("E11").Value = '  -1.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07383'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.173'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.87'
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6628'
$ws.Range("E15").Value = '  -1.84%  '
$ws.Range("D16").Value = '30.332.89'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.47'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007754'
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.461'
$ws.Range("E20").Value = '  +1.51%  '
$ws.Range("D21").Value = '2.148.02'
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '192.15'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.221'
$ws.Range("E24").Value = '  -1.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.457'
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.31'
$ws.Range("E26").Value = '  +1.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.26'
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.940'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.262'
$ws.Range("E30").Value = '  -2.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09170'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.058'
$ws.Range("E32").Value = '  -0.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05076'
$ws.Range("E33").Value = '  -3.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7326'
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.145'
$ws.Range("E35").Value = '  +2.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.714'
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.649'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9194'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.089'
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.922'
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4339'
$ws.Range("E42").Value = '  -3.66%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.14'
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1382'
$ws.Range("E45").Value = '  -1.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.689'
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.585'
$ws.Range("E47").Value = '  +9.01%  '
$ws.Range("E48").Value = '  -9.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.975'
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.34'
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05775'
$ws.Range("E51").Value = '  -1.97%  '
